$wb = $excel.ActiveWorkbook

# --- Demand sheet: move the "t" header from I3 (data row) up into I1 (header row) ---
$wsDemand = $wb.Worksheets.Item("Demand")
$wsDemand.Range("I1").Value = "t"
$wsDemand.Range("I3").ClearContents()

# --- Time_Series_Storage sheet: move the "t" header from A2 (data row) up into A1 (header row) ---
$wsStorage = $wb.Worksheets.Item("Time_Series_Storage")
$wsStorage.Activate()
$wsStorage.Range("A1").Value = "t"
$wsStorage.Range("A2").ClearContents()
$wsStorage.Range("A2").Select()

# --- Make Demand the active sheet/tab (was Time_Series_Storage before) ---
$wsDemand.Activate()
